$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (column C) date bumps from 2023-09-20 (45189) to 2023-09-21 (45190)
#    for every existing data row (2..515).
$ws.Range("C2:C515").Value = 45190

# 2) Row 515 becomes an explicit custom-height row (ht="15" customHeight="1"),
#    matching the formatting already used by the surrounding rows.
$ws.Rows.Item(515).RowHeight = 15

# 3) A brand-new record is appended as row 516.
$ws.Cells.Item(516, 1).Value = "A 44327-2023"

$ws.Cells.Item(516, 2).Value = 45188
$ws.Cells.Item(516, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(516, 3).Value = 45190
$ws.Cells.Item(516, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(516, 4).Value = "JÄMTLANDS LÄN"
$ws.Cells.Item(516, 5).Value = "ÖSTERSUND"

$ws.Cells.Item(516, 7).Value = 1.5
$ws.Cells.Item(516, 8).Value = 0
$ws.Cells.Item(516, 9).Value = 0
$ws.Cells.Item(516, 10).Value = 0
$ws.Cells.Item(516, 11).Value = 0
$ws.Cells.Item(516, 12).Value = 0
$ws.Cells.Item(516, 13).Value = 0
$ws.Cells.Item(516, 14).Value = 0
$ws.Cells.Item(516, 15).Value = 0
$ws.Cells.Item(516, 16).Value = 0
$ws.Cells.Item(516, 17).Value = 0

# Column R keeps the same wrap-text styling used throughout the sheet,
# with no content (same as the other rows' R cell).
$ws.Cells.Item(516, 18).WrapText = $true
